# Auto-generated Excel COM-interop script
# Moves current 'New' sheet rows into 'Previously added', then populates 'New' with fresh scraped rows.

$wb = $excel.ActiveWorkbook
$wsNew = $wb.Worksheets.Item("New")
$wsPrev = $wb.Worksheets.Item("Previously added")

$xlPasteFormats = -4122

# --- Stash style templates from the 'New' sheet (link / plain / date columns) ---
# before we overwrite its data rows, so every new row can be re-stamped with the
# same look (custom hyperlink font, plain data font, custom date format).
$wsNew.Range("A2").Copy()
$wsNew.Range("J1").PasteSpecial($xlPasteFormats)
$wsNew.Range("B2").Copy()
$wsNew.Range("J2").PasteSpecial($xlPasteFormats)
$wsNew.Range("F2").Copy()
$wsNew.Range("J3").PasteSpecial($xlPasteFormats)

# --- Stash style templates from the 'Previously added' sheet (link / plain columns) ---
$wsPrev.Range("A2666").Copy()
$wsPrev.Range("J1").PasteSpecial($xlPasteFormats)
$wsPrev.Range("B2666").Copy()
$wsPrev.Range("J2").PasteSpecial($xlPasteFormats)

# --- Part 1: archive the previous 'New' rows (2-16) onto the bottom of 'Previously added' ---
$oldRows = @(
  @('https://www.ss.com/msg/lv/real-estate/plots-and-lands/aizkraukle-and-reg/kokneses-pag/hndpl.html', '26 000 € (1.50 €/m²)', 'Aizkraukle un raj.', '1.73 ha.', '32600120088', 45945.87777777778),
  @('https://www.ss.com/msg/lv/real-estate/plots-and-lands/bauska-and-reg/vecumnieku-pag/bdnlb.html', '82 500 € (0.75 €/m²)', 'Bauska un raj.', '11 ha.', '32620060046', 45946.44861111111),
  @('https://www.ss.com/msg/lv/real-estate/plots-and-lands/cesis-and-reg/nitaures-pag/inkic.html', '13 000 € (1.30 €/m²)', 'Cēsis un raj.', '1 ha.', '42680060126', 45945.87430555555),
  @('https://www.ss.com/msg/lv/real-estate/plots-and-lands/daugavpils-and-reg/kalkunes-pag/akxpb.html', '9 800 € (0.55 €/m²)', 'Daugavpils un raj.', '17900 m²', '4600042151', 45946.45625),
  @('https://www.ss.com/msg/lv/real-estate/plots-and-lands/daugavpils-and-reg/naujenes-pag/kdcxk.html', '91 500 € (0.74 €/m²)', 'Daugavpils un raj.', '12.35 ha.', '', 45945.975),
  @('https://www.ss.com/msg/lv/real-estate/plots-and-lands/daugavpils-and-reg/laucesas-pag/jngmc.html', '12 000 € (0.40 €/m²)', 'Daugavpils un raj.', '3 ha.', '', 45945.79236111111),
  @('https://www.ss.com/msg/lv/real-estate/plots-and-lands/liepaja-and-reg/pavilosta/eejlm.html', '175 000 € (174.30 €/m²)', 'Liepāja un raj.', '1004 m²', '', 45945.768055555556),
  @('https://www.ss.com/msg/lv/real-estate/plots-and-lands/ogre-and-reg/jumpravas-pag/hbnmx.html', '20 000 € (2.63 €/m²)', 'Ogre un raj.', '7600 m²', '74480020326', 45946.554861111115),
  @('https://www.ss.com/msg/lv/real-estate/plots-and-lands/ogre-and-reg/keguma-l-t/blxmje.html', '35 000 € (1.52 €/m²)', 'Ogre un raj.', '23000 m²', '74290050003', 45946.44930555555),
  @('https://www.ss.com/msg/lv/real-estate/plots-and-lands/preili-and-reg/preili/bfiei.html', '13 800 € (21.50 €/m²)', 'Preiļi un raj.', '642 m²', '', 45946.3375),
  @('https://www.ss.com/msg/lv/real-estate/plots-and-lands/rezekne-and-reg/nautrenu-pag/icgxe.html', '520 200 € (0.45 €/m²)', 'Rēzekne un raj.', '115.60 ha.', '68760010050', 45945.70486111111),
  @('https://www.ss.com/msg/lv/real-estate/plots-and-lands/talsi-and-reg/gibulu-pag/akljb.html', '70 000 € (2.33 €/m²)', 'Talsi un raj.', '3 ha.', '88540160060', 45946.60833333334),
  @('https://www.ss.com/msg/lv/real-estate/plots-and-lands/tukums-and-reg/jaunpils-pag/apopp.html', '140 000 € (0.93 €/m²)', 'Tukums un raj.', '15 ha.', '90560030123', 45945.839583333334),
  @('https://www.ss.com/msg/lv/real-estate/plots-and-lands/valmiera-and-reg/valmiera/hkhpm.html', '42 000 € (30.59 €/m²)', 'Valmiera un raj.', '1373 m²', '', 45945.64444444445),
  @('https://www.ss.com/msg/lv/real-estate/plots-and-lands/ventspils-and-reg/ventspils/bjhge.html', '35 000 € (14.49 €/m²)', 'Ventspils un raj.', '2415 m²', '', 45946.56805555556)
)

$destRow = $wsPrev.UsedRange.Row + $wsPrev.UsedRange.Rows.Count
foreach ($row in $oldRows) {
  $linkCell = $wsPrev.Cells.Item($destRow, 1)
  $linkCell.Value = $row[0]
  $wsPrev.Hyperlinks.Add($linkCell, $row[0])
  $wsPrev.Range("J1").Copy()
  $linkCell.PasteSpecial($xlPasteFormats)

  $wsPrev.Cells.Item($destRow, 2).Value = $row[1]
  $wsPrev.Cells.Item($destRow, 3).Value = $row[2]
  $wsPrev.Cells.Item($destRow, 4).Value = $row[3]
  $wsPrev.Cells.Item($destRow, 5).Value = $row[4]
  $wsPrev.Cells.Item($destRow, 6).Value = [double]$row[5]
  $wsPrev.Range("J2:J2").Copy()
  $wsPrev.Range($wsPrev.Cells.Item($destRow,2), $wsPrev.Cells.Item($destRow,6)).PasteSpecial($xlPasteFormats)

  $destRow = $destRow + 1
}

# --- Part 2: clear the old rows out of 'New' and load this run's freshly scraped rows ---
$wsNew.Range("A2:F16").Clear()

$newRows = @(
  @('https://www.ss.com/msg/lv/real-estate/plots-and-lands/bauska-and-reg/bauska/bldfn.html', '9 000 € (12.68 €/m²)', 'Bauska un raj.', '710 m²', '40500010114', 45946.691666666666),
  @('https://www.ss.com/msg/lv/real-estate/plots-and-lands/cesis-and-reg/marsnenu-pag/jxcmd.html', '6 000 € (0.57 €/m²)', 'Cēsis un raj.', '1.05 ha.', '42640030222', 45947.61319444445),
  @('https://www.ss.com/msg/lv/real-estate/plots-and-lands/cesis-and-reg/amatas-pag/ogbid.html', '180 000 € (1.06 €/m²)', 'Cēsis un raj.', '17 ha.', '', 45947.37986111111),
  @('https://www.ss.com/msg/lv/real-estate/plots-and-lands/dobele-and-reg/dobeles-pag/hojkc.html', '16 900 € (14.02 €/m²)', 'Dobele un raj.', '1205 m²', '', 45946.648611111115),
  @('https://www.ss.com/msg/lv/real-estate/plots-and-lands/gulbene-and-reg/galgauskas-pag/kepbd.html', '62 000 € (0.30 €/m²)', 'Gulbene un raj.', '21 ha.', '', 45947.507638888885),
  @('https://www.ss.com/msg/lv/real-estate/plots-and-lands/jelgava-and-reg/jelgava/bddpgj.html', '89 000 € (24.75 €/m²)', 'Jelgava un raj.', '3596 m²', '09000070224', 45947.63680555555),
  @('https://www.ss.com/msg/lv/real-estate/plots-and-lands/jelgava-and-reg/jelgava/dijlf.html', '40 120 € (40 €/m²)', 'Jelgava un raj.', '1003 m²', '09000210930', 45947.63333333333),
  @('https://www.ss.com/msg/lv/real-estate/plots-and-lands/jelgava-and-reg/jelgava/bdcif.html', '32 000 € (30.19 €/m²)', 'Jelgava un raj.', '1060 m²', '09000100152', 45947.63333333333),
  @('https://www.ss.com/msg/lv/real-estate/plots-and-lands/jelgava-and-reg/cenu-pag/ancip.html', '20 000 € (1.87 €/m²)', 'Jelgava un raj.', '10700 m²', '', 45947.586111111115),
  @('https://www.ss.com/msg/lv/real-estate/plots-and-lands/jelgava-and-reg/jelgava/bdomoc.html', '28 500 € (18.58 €/m²)', 'Jelgava un raj.', '1534 m²', '', 45947.41458333333),
  @('https://www.ss.com/msg/lv/real-estate/plots-and-lands/kraslava-and-reg/dagdas-pag/beddh.html', '20 000 € (0.29 €/m²)', 'Krāslava un raj.', '7 ha.', '60540040296', 45946.75555555556),
  @('https://www.ss.com/msg/lv/real-estate/plots-and-lands/limbadzi-and-reg/salacgriva/bdioj.html', '390 990 € (30 €/m²)', 'Limbaži un raj.', '13033 m²', '', 45947.44236111111),
  @('https://www.ss.com/msg/lv/real-estate/plots-and-lands/madona-and-reg/vestienas-pag/decik.html', '98 761 € (3.29 €/m²)', 'Madona un raj.', '3 ha.', '70960040091', 45946.975),
  @('https://www.ss.com/msg/lv/real-estate/plots-and-lands/ogre-and-reg/ikskile/afmgf.html', '47 000 € (10.56 €/m²)', 'Ogre un raj.', '4450 m²', '', 45947.57916666666),
  @('https://www.ss.com/msg/lv/real-estate/plots-and-lands/ogre-and-reg/kegums/cpfig.html', '40 000 € (4 €/m²)', 'Ogre un raj.', '1 ha.', '74090030087', 45946.85902777778),
  @('https://www.ss.com/msg/lv/real-estate/plots-and-lands/valmiera-and-reg/kocenu-pag/cinix.html', '9 500 € (4.10 €/m²)', 'Valmiera un raj.', '2315 m²', '96640140072', 45947.600694444445),
  @('https://www.ss.com/msg/lv/real-estate/plots-and-lands/valmiera-and-reg/valmiera/bdlhnp.html', '35 000 € (20.72 €/m²)', 'Valmiera un raj.', '1689 m²', '', 45947.495833333334),
  @('https://www.ss.com/msg/lv/real-estate/plots-and-lands/ventspils-and-reg/varves-pag/gjkic.html', '16 500 € (6.60 €/m²)', 'Ventspils un raj.', '2500 m²', '98840010213', 45946.77569444444),
  @('https://www.ss.com/msg/lv/real-estate/plots-and-lands/other/bjmkl.html', '654 321 € (93.47 €/m²)', 'Cits', '7000 m²', '7352036Uk9375S0001Dj', 45946.96736111111)
)

$destRow = 2
foreach ($row in $newRows) {
  $linkCell = $wsNew.Cells.Item($destRow, 1)
  $linkCell.Value = $row[0]
  $wsNew.Hyperlinks.Add($linkCell, $row[0])
  $wsNew.Range("J1").Copy()
  $linkCell.PasteSpecial($xlPasteFormats)

  $wsNew.Cells.Item($destRow, 2).Value = $row[1]
  $wsNew.Cells.Item($destRow, 3).Value = $row[2]
  $wsNew.Cells.Item($destRow, 4).Value = $row[3]
  $wsNew.Cells.Item($destRow, 5).Value = $row[4]
  $wsNew.Range("J2").Copy()
  $wsNew.Range($wsNew.Cells.Item($destRow,2), $wsNew.Cells.Item($destRow,5)).PasteSpecial($xlPasteFormats)

  $wsNew.Cells.Item($destRow, 6).Value = [double]$row[5]
  $wsNew.Range("J3").Copy()
  $wsNew.Cells.Item($destRow, 6).PasteSpecial($xlPasteFormats)

  $destRow = $destRow + 1
}

# --- tidy up: drop the scratch style templates we stashed in column J ---
$wsNew.Range("J1:J3").Clear()
$wsPrev.Range("J1:J2").Clear()

Write-Host "done"
